$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) hold text values (locale-formatted numbers using "." as
# thousands separator), so force text format before assignment to avoid Excel
# auto-converting them into floating point numbers, then restore default style.

$dCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns for rows with refreshed data
$ws.Range("D2").Value = "28.271.04"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "1.816.66"
$ws.Range("E3").Value = "  +4.15%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "327.80"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "0.4334"
$ws.Range("E7").Value = "  +2.66%  "

$ws.Range("D8").Value = "0.3660"
$ws.Range("E8").Value = "  +2.22%  "

$ws.Range("D9").Value = "44.88"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("D10").Value = "0.07676"
$ws.Range("E10").Value = "  +3.42%  "

$ws.Range("D11").Value = "1.142"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "22.02"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").Value = "6.291"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").Value = "7.524"
$ws.Range("E15").Value = "  +4.74%  "

$ws.Range("D16").Value = "1.829.74"
$ws.Range("E16").Value = "  +5.04%  "

$ws.Range("D17").Value = "93.50"
$ws.Range("E17").Value = "  +6.76%  "

$ws.Range("D18").Value = "0.00001081"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").Value = "0.06562"
$ws.Range("E19").Value = "  +6.19%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "17.49"
$ws.Range("E21").Value = "  +3.72%  "

$ws.Range("D22").Value = "6.259"
$ws.Range("E22").Value = "  +2.58%  "

$ws.Range("D23").Value = "28.305.63"
$ws.Range("E23").Value = "  +2.88%  "

$ws.Range("D24").Value = "11.58"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").Value = "2.058"
$ws.Range("E25").Value = "  -11.55%  "

$ws.Range("D26").Value = "162.69"
$ws.Range("E26").Value = "  +6.92%  "

$ws.Range("D27").Value = "20.65"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").Value = "2.040.34"
$ws.Range("E28").Value = "  +5.19%  "

$ws.Range("D29").Value = "2.289"
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("D30").Value = "128.77"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("D31").Value = "1.227"
$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Value = "5.949"
$ws.Range("E32").Value = "  +4.62%  "

$ws.Range("D33").Value = "0.09174"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").Value = "3.486"
$ws.Range("E34").Value = "  -5.55%  "

$ws.Range("D35").Value = "12.97"
$ws.Range("E35").Value = "  +2.43%  "

$ws.Range("D36").Value = "0.02348"
$ws.Range("E36").Value = "  +2.54%  "

$ws.Range("D37").Value = "0.2173"
$ws.Range("E37").Value = "  +2.27%  "

$ws.Range("D38").Value = "5.193"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").Value = "0.6562"
$ws.Range("E39").Value = "  +2.65%  "

$ws.Range("D40").Value = "0.06195"
$ws.Range("E40").Value = "  +1.87%  "

$ws.Range("D41").Value = "1.193"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").Value = "8.123"
$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("D43").Value = "1.436"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D45").Value = "13.88"
$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("D46").Value = "0.6102"
$ws.Range("E46").Value = "  +3.93%  "

$ws.Range("D47").Value = "3.752"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  +2.99%  "

$ws.Range("D51").Value = "0.07006"
$ws.Range("E51").Value = "  +2.30%  "

# Rows 48 and 49 swapped ranking order (Quant now ranks above NEARProtocol)
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "125.69"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.019"
$ws.Range("E49").Value = "  +3.50%  "

# Restore default (Normal) style on the Price cells so only their text content
# differs from the original workbook formatting.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
